# Scheduled-runner refresh of market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, cols H:N)
# across the per-class Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1243
$ws.Range("I15").Value = 1243
$ws.Range("K15").Value = 3729
$ws.Range("M15").Value = -3560

$ws.Range("H28").Value = 692.3077
$ws.Range("I28").Value = 349.875
$ws.Range("J28").Value = 1240.2
$ws.Range("K28").Value = 349.875
$ws.Range("L28").Value = 1240.2
$ws.Range("M28").Value = 135.125
$ws.Range("N28").Value = -2210.2

$ws.Range("H46").Value = 1024.4667
$ws.Range("J46").Value = 1024.4667
$ws.Range("L46").Value = 3073.4001
$ws.Range("N46").Value = -3311.4001

$ws.Range("H60").Value = 1024.4667
$ws.Range("J60").Value = 1024.4667
$ws.Range("L60").Value = 3073.4001
$ws.Range("N60").Value = -4041.4001

$ws.Range("H138").Value = 2826
$ws.Range("I138").Value = 2547.9167
$ws.Range("K138").Value = 7643.750100000001
$ws.Range("M138").Value = -2503.750100000001

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 91784
$ws.Range("I45").Value = 250890
$ws.Range("J45").Value = 866.2857
$ws.Range("K45").Value = 250890
$ws.Range("L45").Value = 866.2857
$ws.Range("M45").Value = -250513
$ws.Range("N45").Value = -1620.2857

$ws.Range("H61").Value = 2067.4
$ws.Range("I61").Value = 2000.7858
$ws.Range("K61").Value = 2000.7858
$ws.Range("M61").Value = -1788.7858

$ws.Range("H97").Value = 60435.883
$ws.Range("I97").Value = 67619.92999999999
$ws.Range("K97").Value = 67619.92999999999
$ws.Range("M97").Value = -67123.92999999999

$ws.Range("H110").Value = 77078616
$ws.Range("I110").Value = 100201050
$ws.Range("J110").Value = 3829.3333
$ws.Range("K110").Value = 100201050
$ws.Range("L110").Value = 3829.3333
$ws.Range("M110").Value = -100199005
$ws.Range("N110").Value = -7919.3333

$ws.Range("H122").Value = 1461.8846
$ws.Range("I122").Value = 1454.1578
$ws.Range("J122").Value = 1482.8572
$ws.Range("K122").Value = 4362.4734
$ws.Range("L122").Value = 4448.571599999999
$ws.Range("M122").Value = -1912.4734
$ws.Range("N122").Value = -9348.571599999999

$ws.Range("H132").Value = 4859.913
$ws.Range("I132").Value = 5171.25
$ws.Range("K132").Value = 15513.75
$ws.Range("M132").Value = -12983.75

$ws.Range("H136").Value = 2067.4
$ws.Range("I136").Value = 2000.7858
$ws.Range("K136").Value = 6002.357400000001
$ws.Range("M136").Value = -3452.357400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2558.5217
$ws.Range("I134").Value = 2824.0557
$ws.Range("J134").Value = 1602.6
$ws.Range("K134").Value = 8472.167099999999
$ws.Range("L134").Value = 4807.799999999999
$ws.Range("M134").Value = -5937.167099999999
$ws.Range("N134").Value = -9877.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21277.355
$ws.Range("I31").Value = 1541.6818
$ws.Range("J31").Value = 29790.785
$ws.Range("K31").Value = 1541.6818
$ws.Range("L31").Value = 29790.785
$ws.Range("M31").Value = -1246.6818
$ws.Range("N31").Value = -30380.785

$ws.Range("H34").Value = 21277.355
$ws.Range("I34").Value = 1541.6818
$ws.Range("J34").Value = 29790.785
$ws.Range("K34").Value = 1541.6818
$ws.Range("L34").Value = 29790.785
$ws.Range("M34").Value = -1339.6818
$ws.Range("N34").Value = -30194.785

$ws.Range("H58").Value = 1948.8182
$ws.Range("I58").Value = 2239
$ws.Range("K58").Value = 2239
$ws.Range("M58").Value = -2036

$ws.Range("H132").Value = 125005944
$ws.Range("I132").Value = 250009650
$ws.Range("J132").Value = 62504090
$ws.Range("K132").Value = 750028950
$ws.Range("L132").Value = 187512270
$ws.Range("M132").Value = -750026420
$ws.Range("N132").Value = -187517330

$ws.Range("H136").Value = 1948.8182
$ws.Range("I136").Value = 2239
$ws.Range("K136").Value = 6717
$ws.Range("M136").Value = -4167

$ws.Range("H137").Value = 63266.668
$ws.Range("J137").Value = 63266.668
$ws.Range("L137").Value = 63266.668
$ws.Range("N137").Value = -73466.66800000001

$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4116
$ws.Range("I56").Value = 4116
$ws.Range("K56").Value = 4116
$ws.Range("M56").Value = -3586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H113").Value = 1649.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1649.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1649.5
$ws.Range("N113").Value = -5989.5
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 985.8889
$ws.Range("I122").Value = 922.36365
$ws.Range("J122").Value = 1085.7142
$ws.Range("K122").Value = 2767.09095
$ws.Range("L122").Value = 3257.1426
$ws.Range("M122").Value = -317.0909499999998
$ws.Range("N122").Value = -8157.142599999999

$ws.Range("H132").Value = 2849.6
$ws.Range("I132").Value = 2633.0527
$ws.Range("J132").Value = 3535.3333
$ws.Range("K132").Value = 7899.158100000001
$ws.Range("L132").Value = 10605.9999
$ws.Range("M132").Value = -5369.158100000001
$ws.Range("N132").Value = -15665.9999

$ws.Range("H140").Value = 124754.29
$ws.Range("J140").Value = 124754.29
$ws.Range("L140").Value = 124754.29
$ws.Range("N140").Value = -135114.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1934.25
$ws.Range("I61").Value = 1878.8889
$ws.Range("J61").Value = 2005.4286
$ws.Range("K61").Value = 1878.8889
$ws.Range("L61").Value = 2005.4286
$ws.Range("M61").Value = -1676.8889
$ws.Range("N61").Value = -2409.4286

$ws.Range("H113").Value = 1934.25
$ws.Range("I113").Value = 1878.8889
$ws.Range("J113").Value = 2005.4286
$ws.Range("K113").Value = 1878.8889
$ws.Range("L113").Value = 2005.4286
$ws.Range("M113").Value = 291.1111000000001
$ws.Range("N113").Value = -6345.4286

$ws.Range("H132").Value = 4387.2085
$ws.Range("I132").Value = 4252.143
$ws.Range("K132").Value = 12756.429
$ws.Range("M132").Value = -10226.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 998.5
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 998.125
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 2994.375
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7334.375

$ws.Range("H132").Value = 2268.25
$ws.Range("I132").Value = 2245.394
$ws.Range("J132").Value = 2376
$ws.Range("K132").Value = 6736.181999999999
$ws.Range("L132").Value = 7128
$ws.Range("M132").Value = -4206.181999999999
$ws.Range("N132").Value = -12188

$ws.Range("H136").Value = 2055.0952
$ws.Range("I136").Value = 564
$ws.Range("K136").Value = 1692
$ws.Range("M136").Value = 858
